$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 / Row 25 swap (Dai <-> Cosmos) plus value updates ---
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"

$ws.Range("D2").Value = "42.536.25"
$ws.Range("E2").Value = "  +1.64%  "

$ws.Range("D3").Value = "2.254.10"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.46%  "

$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.99"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +10.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.28"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.38%  "

$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").Value = "2.265.98"
$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("D17").Value = "42.388.09"
$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("E18").Value = "  +4.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.69%  "

$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +25.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.57"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("E26").Value = "  -2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.32"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.27%  "

$ws.Range("E28").Value = "  +1.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.64"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.73"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0827"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.33%  "

$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.31"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.82%  "

$ws.Range("E35").Value = "  +0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.57"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0315"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.19"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.83"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.98"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.17%  "

$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "108.26"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.82"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.59%  "

$ws.Range("E45").Value = "  +3.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("E47").Value = "  +5.73%  "

$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.15"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("E51").Value = "  +0.82%  "
